$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ConDA (NewsCLIPpings)")

$ws.Range("F4").Value = 0.72160000000000002
$ws.Range("G4").Value = 0.77980000000000005
$ws.Range("H4").Value = 0.83730000000000004
$ws.Range("I4").Value = 0.80269999999999997
